$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of worksheet row -> new dSF (column F) value
$changes = @{
    2  = 0
    8  = 0
    13 = 1
    14 = 4
    15 = 2
    16 = 0
    19 = 0
    35 = 1
    38 = -1
    40 = -2
    44 = 0
    46 = 1
    52 = 0
    54 = 0
    58 = -4
    60 = -2
    61 = -1
    63 = -1
    64 = 1
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
